$wb = $excel.ActiveWorkbook

# --- Sheet "Accuracy" (sheet1) ---
$ws1 = $wb.Worksheets.Item("Accuracy")

$ws1.Range("F1").Value = "TDD Given Test"
$ws1.Range("F5").Value = 67.599999999999994
$ws1.Range("F6").Value = 74.400000000000006
$ws1.Range("F5:F6").NumberFormat = "0.00"

$ws1.Range("H15").Select()

# --- Sheet "Errors" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Errors")

$ws2.Range("N1:P1").Merge()
$ws2.Range("N1").Value = "TDD Given Test"
$ws2.Range("N1:P1").HorizontalAlignment = -4108

$ws2.Range("N2").Value = "Assertion"
$ws2.Range("O2").Value = "Runtime"
$ws2.Range("P2").Value = "Compilation"

$ws2.Range("N6").Value = 149
$ws2.Range("O6").Value = 9
$ws2.Range("P6").Value = 4

$ws2.Range("N7").Value = 119
$ws2.Range("O7").Value = 9
$ws2.Range("P7").Value = 0

$ws2.Columns.Item(16).ColumnWidth = 11

$ws2.Range("P12").Select()
